$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.031.03"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = "'1.639.74"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'215.08"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').Value = "'0.5055"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('D7').Value = "'1.008"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'0.2572"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('D11').Value = "'0.07725"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('D12').Value = "'1.647.12"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = "'4.245"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = "'1.865.95"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = "'0.5445"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.58%  '
$ws.Range('D16').Value = "'0.0₅7904"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('D17').Value = "'63.63"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = "'26.007.41"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').Value = "'1.008"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = "'203.84"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('D22').Value = "'9.975"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').Value = "'5.949"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = "'1.928"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.44%  '
$ws.Range('D26').Value = "'141.38"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('D27').Value = "'0.1152"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('D28').Value = "'15.71"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').Value = "'6.726"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.48%  '
$ws.Range('D30').Value = "'0.05048"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.29%  '
$ws.Range('D31').Value = "'1.239"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('D32').Value = "'3.250"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('D33').Value = "'3.193"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').Value = "'1.541"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').Value = "'2.337"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.17%  '
$ws.Range('D36').Value = "'2.631"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.76%  '
$ws.Range('D37').Value = "'0.8899"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.93%  '
$ws.Range('D38').Value = "'0.5619"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').Value = "'1.144.34"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('D40').Value = "'0.01572"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D41').Value = "'2.564"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('D42').Value = "'1.008"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = "'5.654"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = "'0.8093"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.37%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = "'1.778.15"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').Value = "'0.4526"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = "'54.94"
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = "'0.05035"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.04%  '
